# time_log.xlsx edit: "Finish clean up the words"
#
# Row 59 (2021-11-24 entry, JS101 course): the Hours value was bumped from
# 1.25 to 1.5, and the note text was changed from the one-off
# "Finished 1 small problem, worked on a second" to the already-used
# "Finished 2 small problems" (matching the wording used elsewhere, e.g.
# row 43). Excel's shared-string table drops the now-unreferenced string
# automatically on save.
#
# Downstream weekly/grand totals (D64 SUM and C65 SUBTOTAL) recalc on their
# own since they're live formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C59").Value = 1.5
$ws.Range("D59").Value = "Finished 2 small problems"

$ws.Calculate()
